# Checking for bad header
#
# Adds a new "BadHeader" worksheet that mirrors the "WithTitle" table but
# uses mismatched header captions ("First"/"Last" instead of
# "First name"/"Last name") so the ORM's header-detection logic can be
# exercised against a table whose header row doesn't match the expected
# column names. Also updates the selections that were left behind on the
# other two sheets when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# --- "Tab" sheet: just a different cell was selected when saving -------
$ws1 = $wb.Worksheets.Item("Tab")
$ws1.Range("C8").Select()

# --- "WithTitle" sheet: the little data table is now selected ----------
$ws2 = $wb.Worksheets.Item("WithTitle")
$ws2.Activate()
$ws2.Range("B7:D11").Select()

# --- New "BadHeader" sheet, appended after the existing sheets ---------
$srcTable = $ws2.Range("B7:D11")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$badHeader = $wb.Worksheets.Add($null, $lastSheet)
$badHeader.Name = "BadHeader"

# Copy the table (values + styles) from "WithTitle" into A1:C5
$destTable = $badHeader.Range("A1:C5")
$srcTable.Copy($destTable)

# Overwrite the header row with intentionally "bad" captions
$badHeader.Range("A1").Value = "First"
$badHeader.Range("B1").Value = "Last"

# Restore the taller header/row-3 heights that come with the source rows
$badHeader.Rows.Item(1).RowHeight = 25.5
$badHeader.Rows.Item(3).RowHeight = 25.5

# Leave the new sheet active, with C1 selected
$badHeader.Activate()
$badHeader.Range("C1").Select()
